$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new test case rows (data & run-mode) below the existing ones.
$ws.Range("A5").Value = "testSearchByOrg"
$ws.Range("C5").Value = "Y"

$ws.Range("A6").Value = "testSearchByOrgKeyword"
$ws.Range("C6").Value = "Y"

# Move the active selection to the last added cell, matching the authored workbook.
$ws.Range("C6").Select() | Out-Null
